$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 15)

$data = @(
    @(9.970763543028557, 8.528491046257558, 13.59532827535271, 36.94469657144347, 3.688423916727922, 28.75196668293579, 10.3854994726774, 19.26589311391507, 10.41771353776731, 28.33118522513302),
    @(9.901449391236348, 8.472637652341685, 13.58267831468798, 37.12165219285041, 3.690703559354318, 28.89341137294942, 10.40696812928699, 18.66807038968713, 10.43272694187783, 28.49810661952498),
    @(9.859922379524233, 8.438888176939864, 13.57668352587949, 37.24012055144738, 3.692176221356632, 28.98724341742677, 10.42123783384365, 18.29119183074582, 10.44298308342464, 28.60768619690392),
    @(9.843272382012442, 8.425281244204822, 13.57468859291661, 37.29085823204711, 3.692794746457167, 29.02723321562139, 10.42732678634423, 18.13534225733512, 10.44742392753713, 28.65412057357878),
    @(9.840524499956006, 8.423030923338736, 13.57438445743427, 37.29943159006292, 3.69289856533138, 29.03397922555163, 10.42835440822706, 18.10933268645005, 10.44817712496356, 28.66193838212665),
    @(9.859696711157753, 8.43870406481798, 13.57665480468494, 37.24079486383425, 3.692184488409221, 28.98777564316702, 10.42131884194351, 18.28909889872783, 10.44304191533787, 28.60830522445789),
    @(9.946657750836181, 8.509125362261701, 13.59060016904011, 37.00366877871611, 3.689194831842525, 28.79928489795135, 10.3926763086637, 19.06191478586019, 10.42267500037706, 28.38726665045428),
    @(10.12478736512774, 8.651128174590188, 13.6319010008117, 36.61689846480569, 3.683908269873577, 28.48524419592491, 10.34512321831443, 20.49164425940205, 10.39095247769416, 28.01018810057037),
    @(10.25948042876953, 8.757273652212133, 13.67059602711632, 36.3808934185132, 3.680371679614317, 28.28864113807828, 10.31541485173778, 21.48036828859746, 10.37263098075792, 27.76769643750288),
    @(10.32140781665694, 8.805833472468303, 13.68997563848075, 36.28408286501337, 3.678837430876387, 28.20666161272857, 10.30303027278246, 21.91509313344502, 10.36537318226339, 27.66492340738206),
    @(10.34493705609971, 8.824250680759555, 13.69756620973748, 36.24894825078015, 3.678267111870657, 28.17669465845869, 10.29850264940616, 22.07742797072271, 10.36277921520776, 27.62709364057483),
    @(10.33986636115899, 8.820283101378665, 13.69592029912707, 36.25644714060907, 3.678389466688351, 28.18310060687137, 10.29947054782865, 22.04256972433711, 10.3633310122714, 27.63519249293137),
    @(10.32334209334767, 8.807348154270077, 13.69059509158616, 36.28116169896251, 3.67879029690282, 28.20417460332064, 10.30265453469941, 21.9284949824232, 10.36515668306384, 27.6617892967688),
    @(10.31323028826693, 8.799428539959216, 13.68736594694568, 36.296498986325, 3.679037204620083, 28.21722342092041, 10.30462592518114, 21.85831982132543, 10.36629505406463, 27.67822244454947),
    @(10.2554451834422, 8.754104747258031, 13.66936497916188, 36.3874332209617, 3.680473442126566, 28.29414913366381, 10.31624692141509, 21.45164441571012, 10.37312692525832, 27.77456493317239),
    @(10.22015262194716, 8.726362566813323, 13.6587746126118, 36.44592641726062, 3.681373586020865, 28.34325362956509, 10.32366518701709, 21.19821633636745, 10.37759353424211, 27.83560147833258),
    @(10.1999158436266, 8.710432404145926, 13.6528508067439, 36.48056260924903, 3.681898346523312, 28.37219876508735, 10.32803834916893, 21.05104093370409, 10.38026398422747, 27.87141730093873),
    @(10.19307523509761, 8.705043585646322, 13.6508739814651, 36.4924600239477, 3.682077229036349, 28.38211943528782, 10.32953730629513, 21.00097159250418, 10.38118557680806, 27.8836656045658),
    @(10.22390321159538, 8.729313113843389, 13.65988466581008, 36.4395969398693, 3.681277037859502, 28.33795373243152, 10.3228644932717, 21.22534107181615, 10.37710756727454, 27.82903060232542),
    @(10.32819365965316, 8.811146769361814, 13.69215242644927, 36.2738609671509, 3.678672274267098, 28.19795539489759, 10.30171492257421, 21.96206444402768, 10.36461625231962, 27.65394760419743),
    @(10.39680490185892, 8.86479186943693, 13.71470792047979, 36.17443821334575, 3.6770320651446, 28.11273766299666, 10.28883742277931, 22.43018400832426, 10.35735223002799, 27.54586537681451),
    @(10.3601496146613, 8.836149177323064, 13.70253667246504, 36.22668532240992, 3.677901806395239, 28.15764385181899, 10.29562403202046, 22.18160003976415, 10.36114699162075, 27.60296892343705),
    @(10.22220740257269, 8.727979110975591, 13.65938229736196, 36.4424553583348, 3.681320664677299, 28.34034759121084, 10.32322614953314, 21.21308256765268, 10.37732695351494, 27.83199903697928),
    @(10.07586771057267, 8.612352239008182, 13.61924938456711, 36.71310671572161, 3.685277137469378, 28.56422527912518, 10.3570676803089, 20.11508051259309, 10.39865703894871, 28.10614427447515)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowData[$j]
    }
}
